# Add a new "case" worksheet (with test-case sample data) as the last tab,
# mirroring what a user would do in Excel: insert a sheet after the last
# existing tab, rename it, populate headers + two sample rows, autosize the
# columns, and leave the selection on the new sheet the way it was left by
# the author (cell G5).

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet so it lands at the end
# of the tab strip (Worksheets.Add defaults to inserting before the active
# sheet, so we explicitly pass the current last sheet as the "After" sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$caseSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$caseSheet.Name = "case"

# Header row.
$caseSheet.Range("A1").Value = "title"
$caseSheet.Range("B1").Value = "state"
$caseSheet.Range("C1").Value = "identifier"
$caseSheet.Range("D1").Value = "tags"
$caseSheet.Range("E1").Value = "description"
$caseSheet.Range("F1").Value = "type"
$caseSheet.Range("G1").Value = "priority"

# First data row (title/identifier/tags/description filled left-to-right,
# "state" for row 3 entered before the rest of that row - matches how the
# original test data was authored).
$caseSheet.Range("A2").Value = "Test Case title - 1"
$caseSheet.Range("B2").Value = "Enquiring"
$caseSheet.Range("C2").Value = "Test identifier - 1"
$caseSheet.Range("D2").Value = "Test tags - 1"
$caseSheet.Range("E2").Value = "Test description - 1"

# Second data row.
$caseSheet.Range("B3").Value = "Reviewing"
$caseSheet.Range("A3").Value = "Test Case title - 2"
$caseSheet.Range("C3").Value = "Test identifier - 2"
$caseSheet.Range("D3").Value = "Test tags - 2"
$caseSheet.Range("E3").Value = "Test description - 2"

# "type" and "priority" columns filled in afterwards.
$caseSheet.Range("F2").Value = "Complaint"
$caseSheet.Range("G2").Value = "High"

$caseSheet.Range("F3").Value = "General Support"
$caseSheet.Range("G3").Value = "Normal"

# Autosize the columns to fit their content (as was done for every other
# sheet in this workbook).
$caseSheet.Columns.Item(1).ColumnWidth = 15.333333333333334
$caseSheet.Columns.Item(2).ColumnWidth = 8.666666666666666
$caseSheet.Columns.Item(3).ColumnWidth = 15.499999999999998
$caseSheet.Columns.Item(4).ColumnWidth = 10.5
$caseSheet.Columns.Item(5).ColumnWidth = 17.166666666666668
$caseSheet.Columns.Item(6).ColumnWidth = 14.666666666666666
$caseSheet.Columns.Item(7).ColumnWidth = 6.666666666666667

# Leave the selection where the author left it.
$caseSheet.Range("G5").Select() | Out-Null
